# Refresh the crypto price snapshot (columns D "Price" and E "Volume(1h)")
# for rows 2-51, matching the Mon Mar 20 11:34:40 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'28.315.55"
$ws.Range("E2").Value2 = "  +4.31%  "
$ws.Range("D3").Value2 = "'1.789.71"
$ws.Range("E3").Value2 = "  +0.49%  "
$ws.Range("E4").Value2 = "  -0.29%  "
$ws.Range("D5").Value2 = "'339.41"
$ws.Range("E5").Value2 = "  +0.61%  "
$ws.Range("D6").Value2 = "'0.9988"
$ws.Range("E6").Value2 = "  -0.21%  "
$ws.Range("D7").Value2 = "'0.3835"
$ws.Range("E7").Value2 = "  -2.67%  "
$ws.Range("D8").Value2 = "'0.3452"
$ws.Range("E8").Value2 = "  +1.00%  "
$ws.Range("D9").Value2 = "'47.12"
$ws.Range("E9").Value2 = "  -1.43%  "
$ws.Range("D10").Value2 = "'1.157"
$ws.Range("E10").Value2 = "  -2.63%  "
$ws.Range("D11").Value2 = "'0.07421"
$ws.Range("E11").Value2 = "  +0.05%  "
$ws.Range("D12").Value2 = "'23.31"
$ws.Range("E12").Value2 = "  +8.26%  "
$ws.Range("E13").Value2 = "  -0.37%  "
$ws.Range("D14").Value2 = "'6.479"
$ws.Range("E14").Value2 = "  +0.83%  "
$ws.Range("D15").Value2 = "'7.362"
$ws.Range("E15").Value2 = "  +3.98%  "
$ws.Range("D16").Value2 = "'1.786.96"
$ws.Range("E16").Value2 = "  +0.42%  "
$ws.Range("D17").Value2 = "'0.00001081"
$ws.Range("E17").Value2 = "  -0.72%  "
$ws.Range("D18").Value2 = "'0.06677"
$ws.Range("E18").Value2 = "  +0.13%  "
$ws.Range("D19").Value2 = "'82.46"
$ws.Range("E19").Value2 = "  -0.99%  "
$ws.Range("D20").Value2 = "'0.9994"
$ws.Range("E20").Value2 = "  -0.15%  "
$ws.Range("D21").Value2 = "'17.55"
$ws.Range("E21").Value2 = "  +0.31%  "
$ws.Range("D22").Value2 = "'6.479"
$ws.Range("E22").Value2 = "  +0.08%  "
$ws.Range("D23").Value2 = "'28.308.24"
$ws.Range("E23").Value2 = "  +4.22%  "
$ws.Range("D24").Value2 = "'12.13"
$ws.Range("E24").Value2 = "  -1.44%  "
$ws.Range("D25").Value2 = "'2.353"
$ws.Range("E25").Value2 = "  -1.01%  "
$ws.Range("D26").Value2 = "'1.448"
$ws.Range("E26").Value2 = "  +0.29%  "
$ws.Range("D27").Value2 = "'20.80"
$ws.Range("E27").Value2 = "  -1.39%  "
$ws.Range("D28").Value2 = "'2.435"
$ws.Range("E28").Value2 = "  -2.31%  "
$ws.Range("D29").Value2 = "'154.90"
$ws.Range("E29").Value2 = "  -0.71%  "
$ws.Range("D30").Value2 = "'135.68"
$ws.Range("E30").Value2 = "  +1.06%  "
$ws.Range("D31").Value2 = "'1.988.13"
$ws.Range("E31").Value2 = "  +0.37%  "
$ws.Range("D32").Value2 = "'6.140"
$ws.Range("E32").Value2 = "  +2.90%  "
$ws.Range("D33").Value2 = "'3.971"
$ws.Range("E33").Value2 = "  +0.13%  "
$ws.Range("D34").Value2 = "'0.08938"
$ws.Range("E34").Value2 = "  +2.69%  "
$ws.Range("D35").Value2 = "'12.84"
$ws.Range("E35").Value2 = "  -0.91%  "
$ws.Range("D36").Value2 = "'0.02434"
$ws.Range("E36").Value2 = "  +1.48%  "
$ws.Range("D37").Value2 = "'0.6890"
$ws.Range("E37").Value2 = "  +1.78%  "
$ws.Range("D38").Value2 = "'5.369"
$ws.Range("E38").Value2 = "  -0.24%  "
$ws.Range("D39").Value2 = "'0.06405"
$ws.Range("E39").Value2 = "  +0.52%  "
$ws.Range("D40").Value2 = "'0.2178"
$ws.Range("E40").Value2 = "  -0.90%  "
$ws.Range("E41").Value2 = "  +0.98%  "
$ws.Range("D42").Value2 = "'1.502"
$ws.Range("E42").Value2 = "  -6.65%  "
$ws.Range("D43").Value2 = "'8.318"
$ws.Range("E43").Value2 = "  -1.09%  "
$ws.Range("D44").Value2 = "'14.25"
$ws.Range("E44").Value2 = "  -0.79%  "
$ws.Range("D45").Value2 = "'0.9985"
$ws.Range("E45").Value2 = "  -0.18%  "
$ws.Range("D46").Value2 = "'0.6332"
$ws.Range("E46").Value2 = "  -0.53%  "
$ws.Range("D47").Value2 = "'3.880"
$ws.Range("E47").Value2 = "  +0.61%  "
$ws.Range("D48").Value2 = "'133.73"
$ws.Range("E48").Value2 = "  +1.64%  "
$ws.Range("D49").Value2 = "'2.091"
$ws.Range("E49").Value2 = "  -1.72%  "
$ws.Range("D50").Value2 = "'0.07494"
$ws.Range("E50").Value2 = "  +5.52%  "
$ws.Range("D51").Value2 = "'1.214"
$ws.Range("E51").Value2 = "  +7.13%  "
